$wb = $excel.ActiveWorkbook

# --- Sheet "OR correctBF by interventn": add IYCF row (row 4) ---
$ws18 = $wb.Worksheets.Item("OR correctBF by interventn")
$ws18.Range("A4").Value = "IYCF"
$ws18.Range("B4").Value = 5.16
$ws18.Range("C4").Value = 5.16
$ws18.Range("D4").Value = 1.82
$ws18.Range("E4").Value = 1.82
$ws18.Range("F4").Value = 1

# --- Sheet "OR stunting by compfeeding": add rows 7-10 (IYCF variants) ---
$ws17 = $wb.Worksheets.Item("OR stunting by compfeeding")
$ws17.Range("A7").Value = "Complementary feeding (food secure with IYCF)"
$ws17.Range("B7").Value = 1
$ws17.Range("C7").Value = 1
$ws17.Range("D7").Value = 1
$ws17.Range("E7").Value = 1
$ws17.Range("F7").Value = 1

$ws17.Range("A8").Value = "Complementary feeding (food secure without IYCF)"
$ws17.Range("B8").Value = 1
$ws17.Range("C8").Value = 1
$ws17.Range("D8").Value = 1.43
$ws17.Range("E8").Value = 1.43
$ws17.Range("F8").Value = 1

$ws17.Range("A9").Value = "Complementary feeding (food insecure with IYCF and supplementation)"
$ws17.Range("B9").Value = 1
$ws17.Range("C9").Value = 1
$ws17.Range("D9").Value = 1.6
$ws17.Range("E9").Value = 1.6
$ws17.Range("F9").Value = 1

$ws17.Range("A10").Value = "Complementary feeding (food insecure with neither IYCF nor supplementation)"
$ws17.Range("B10").Value = 1
$ws17.Range("C10").Value = 1
$ws17.Range("D10").Value = 2.39
$ws17.Range("E10").Value = 2.39
$ws17.Range("F10").Value = 1
